$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "this one"
$ws.Range("B2").Value = "これ"
$ws.Range("A3").Value = "that one"
$ws.Range("B3").Value = "それ"
$ws.Range("A4").Value = "that one (over there)"
$ws.Range("B4").Value = "あれ"
$ws.Range("A5").Value = "which one"
$ws.Range("B5").Value = "どれ"
$ws.Range("A6").Value = "this..."
$ws.Range("B6").Value = "この"
$ws.Range("A7").Value = "that..."
$ws.Range("B7").Value = "その"
$ws.Range("A8").Value = "that... (over there)"
$ws.Range("B8").Value = "あの"
$ws.Range("A9").Value = "which..."
$ws.Range("B9").Value = "どの"
$ws.Range("A10").Value = "here"
$ws.Range("B10").Value = "ここ"
$ws.Range("A11").Value = "there"
$ws.Range("B11").Value = "そこ"
$ws.Range("A12").Value = "over there"
$ws.Range("B12").Value = "あそこ"
$ws.Range("A13").Value = "where"
$ws.Range("B13").Value = "どこ"
$ws.Range("A14").Value = "who"
$ws.Range("B14").Value = "だれ"
$ws.Range("A21").Value = "pencil"
$ws.Range("B21").Value = "えんぴつ"
$ws.Range("A22").Value = "umbrella"
$ws.Range("B22").Value = "かさ"
$ws.Range("A23").Value = "bag"
$ws.Range("B23").Value = "かばん"
$ws.Range("A24").Value = "shoes"
$ws.Range("B24").Value = "くつ"
$ws.Range("A25").Value = "wallet"
$ws.Range("B25").Value = "さいふ"
$ws.Range("A26").Value = "jeans"
$ws.Range("B26").Value = "ジーンズ"
$ws.Range("A27").Value = "dictionary"
$ws.Range("B27").Value = "じしょ"
$ws.Range("A28").Value = "bicycle"
$ws.Range("B28").Value = "じてんしゃ"
$ws.Range("A29").Value = "newspaper"
$ws.Range("B29").Value = "しんぶん"
$ws.Range("A30").Value = "T-shirt"
$ws.Range("B30").Value = "Ｔシャツ"
$ws.Range("A31").Value = "watch; clock"
$ws.Range("B31").Value = "とけい"
$ws.Range("A32").Value = "notebook"
$ws.Range("B32").Value = "ノート"
$ws.Range("A33").Value = "pen"
$ws.Range("B33").Value = "ペン"
$ws.Range("A34").Value = "hat; cap"
$ws.Range("B34").Value = "ぼうし"
$ws.Range("A35").Value = "book"
$ws.Range("B35").Value = "ほん"
$ws.Range("A36").Value = "cafe"
$ws.Range("B36").Value = "きっさてん"
$ws.Range("A37").Value = "bank"
$ws.Range("B37").Value = "ぎんこう"
$ws.Range("A38").Value = "toilet; restroom"
$ws.Range("B38").Value = "トイレ"
$ws.Range("A39").Value = "library"
$ws.Range("B39").Value = "としょかん"
$ws.Range("A40").Value = "post office"
$ws.Range("B40").Value = "ゆうびんきょく"
$ws.Range("A41").Value = "how much"
$ws.Range("B41").Value = "いくら"
$ws.Range("A42").Value = "...yen"
$ws.Range("B42").Value = "～えん"
$ws.Range("A43").Value = "expensive; high"
$ws.Range("B43").Value = "たかい"
$ws.Range("A44").Value = "Welcome (to our store)."
$ws.Range("B44").Value = "いらっしゃいませ"
$ws.Range("A45").Value = "..., please."
$ws.Range("B45").Value = "（～を）おねがいします"
$ws.Range("A46").Value = "Please give me..."
$ws.Range("B46").Value = "（～を）ください"
$ws.Range("A47").Value = "then...; if that is the case,..."
$ws.Range("B47").Value = "じゃあ"
$ws.Range("A48").Value = "Please.; Here it is."
$ws.Range("B48").Value = "どうぞ"
$ws.Range("A49").Value = "Thank you."
$ws.Range("B49").Value = "どうも"
$ws.Range("A50").Value = "set meal; special (of the day)"
$ws.Range("B50").Value = "ていしょく"
$ws.Range("A51").Value = "spaghetti"
$ws.Range("B51").Value = "スパゲッティ"
$ws.Range("A52").Value = "curry"
$ws.Range("B52").Value = "カレー"
$ws.Range("A53").Value = "soba"
$ws.Range("B53").Value = "そば"
$ws.Range("A54").Value = "udon"
$ws.Range("B54").Value = "うどん"
$ws.Range("A55").Value = "sandwich"
$ws.Range("B55").Value = "サンドイッチ"
$ws.Range("A56").Value = "ice cream"
$ws.Range("B56").Value = "アイスクリーム"
$ws.Range("A57").Value = "ramen"
$ws.Range("B57").Value = "ラーメン"
$ws.Range("A58").Value = "hamburger"
$ws.Range("B58").Value = "ハンバーガー"
$ws.Range("A59").Value = "coffee"
$ws.Range("B59").Value = "コーヒー"
$ws.Range("A60").Value = "cola"
$ws.Range("B60").Value = "コーラ"
$ws.Range("A61").Value = "salad"
$ws.Range("B61").Value = "サラダ"
$ws.Range("A62").Value = "black tea"
$ws.Range("B62").Value = "こうちゃ"
$ws.Range("A63").Value = "juice"
$ws.Range("B63").Value = "ジュース"
$ws.Range("A64").Value = "milk"
$ws.Range("B64").Value = "ミルク"
$ws.Range("A65").Value = "blackboard"
$ws.Range("B65").Value = "こくばん"
$ws.Range("A66").Value = "(electric) light"
$ws.Range("B66").Value = "でんき"
$ws.Range("A67").Value = "door (western-style)"
$ws.Range("B67").Value = "ドア"
$ws.Range("A68").Value = "curtain; curtains"
$ws.Range("B68").Value = "カーテン"
$ws.Range("A69").Value = "window"
$ws.Range("B69").Value = "まど"
$ws.Range("A70").Value = "chair; stool"
$ws.Range("B70").Value = "いす"
$ws.Range("A71").Value = "desk"
$ws.Range("B71").Value = "つくえ"
$ws.Range("A72").Value = "eraser"
$ws.Range("B72").Value = "けしゴム"
$ws.Range("A73").Value = "book"
$ws.Range("B73").Value = "ほん"
$ws.Range("A74").Value = "dictionary"
$ws.Range("B74").Value = "じしょ"
$ws.Range("A75").Value = "bag"
$ws.Range("B75").Value = "かばん"
$ws.Range("A76").Value = "pencil"
$ws.Range("B76").Value = "えんぴつ"
$ws.Range("A77").Value = "pen"
$ws.Range("B77").Value = "ペン"
$ws.Range("A78").Value = "Do you understand?"
$ws.Range("B78").Value = "わかりましたか。"
$ws.Range("A79").Value = "I understand./I understood."
$ws.Range("B79").Value = "わかりました。"
$ws.Range("A80").Value = "I don't understand./I don't know."
$ws.Range("B80").Value = "わかりません。"
$ws.Range("A81").Value = "Please say it slowly."
$ws.Range("B81").Value = "ゆっくりいってください。"
$ws.Range("A82").Value = "Please say it again."
$ws.Range("B82").Value = "もういちどいってください。"
$ws.Range("A83").Value = "Please wait."
$ws.Range("B83").Value = "ちょっとまってください。"
$ws.Range("A84").Value = "Please listen./Please ask."
$ws.Range("B84").Value = "きいてください。"
$ws.Range("A85").Value = "Please look at page 10."
$ws.Range("B85").Value = "10ページをみてください。"
